# Baseboard_v0106: mark two CPLD bugs as fixed.
#  D-0020567 --- <EVT FW BB CPLD> Register 80h 0b/1b is reversed to reflect the low active signal.  (row 1)
#  D-0020797 --- <EVT FW BB CPLD> BB CPLD1/2 couldn't be reset by BB_CPLD1/2_RST_L signal.          (row 2)
# Add a new "Fix it in Baseboard_v0106, please verify" note in column H for both rows, styled in red.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$note = "Fix it in Baseboard_v0106, please verify"

$cells = $ws.Range("H1:H2")
$cells.Value = $note
$cells.NumberFormat = "@"
$cells.Font.Color = 255
